$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 794 (shifts old rows 794-835 down to 795-836),
# matching the "dimension" growing from A1:D835 to A1:D836.
$ws.Rows.Item(794).Insert()

# Force column A to be read as plain text (not auto-parsed into a date
# serial) before writing the date string, then populate the new row.
$ws.Range("A794").NumberFormat = "@"
$ws.Range("A794").Value = "2026/02/10"
$ws.Range("B794").Value = "火"
$ws.Range("C794").Value = 18
$ws.Range("D794").Value = 201

# Drop the temporary text-number-format override so the new row's cells
# carry no explicit style, matching their sibling rows.
$ws.Range("A794:D794").ClearFormats()
